$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 379.33334
$ws.Range("I2").Value = 194.25
$ws.Range("K2").Value = 194.25
$ws.Range("M2").Value = -81.25
$ws.Range("H38").Value = 2166.7273
$ws.Range("I38").Value = 975.2857
$ws.Range("J38").Value = 4251.75
$ws.Range("K38").Value = 2925.8571
$ws.Range("L38").Value = 12755.25
$ws.Range("M38").Value = -2553.8571
$ws.Range("N38").Value = -13499.25
$ws.Range("H86").Value = 6584610
$ws.Range("I86").Value = 6841
$ws.Range("J86").Value = 17547558
$ws.Range("K86").Value = 6841
$ws.Range("L86").Value = 17547558
$ws.Range("M86").Value = -5718
$ws.Range("N86").Value = -17549804
$ws.Range("H89").Value = 6584610
$ws.Range("I89").Value = 6841
$ws.Range("J89").Value = 17547558
$ws.Range("K89").Value = 34205
$ws.Range("L89").Value = 87737790
$ws.Range("M89").Value = -28589
$ws.Range("N89").Value = -87749022
$ws.Range("H135").Value = 15156987
$ws.Range("I135").Value = 20834898
$ws.Range("K135").Value = 187514082
$ws.Range("M135").Value = -187511547

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 725
$ws.Range("I30").Value = 750
$ws.Range("J30").Value = 700
$ws.Range("K30").Value = 750
$ws.Range("L30").Value = 700
$ws.Range("M30").Value = -600
$ws.Range("N30").Value = -1000
$ws.Range("H31").Value = 17516.5
$ws.Range("I31").Value = 15018.857
$ws.Range("K31").Value = 15018.857
$ws.Range("M31").Value = -14724.857
$ws.Range("H61").Value = 5689.5
$ws.Range("I61").Value = 4730.273
$ws.Range("K61").Value = 4730.273
$ws.Range("M61").Value = -4518.273
$ws.Range("H74").Value = 1771.579
$ws.Range("J74").Value = 3862
$ws.Range("L74").Value = 3862
$ws.Range("N74").Value = -5610
$ws.Range("H77").Value = 1771.579
$ws.Range("J77").Value = 3862
$ws.Range("L77").Value = 19310
$ws.Range("N77").Value = -28046
$ws.Range("H88").Value = 4637189.5
$ws.Range("I88").Value = 19115.166
$ws.Range("J88").Value = 6946227
$ws.Range("K88").Value = 19115.166
$ws.Range("L88").Value = 6946227
$ws.Range("M88").Value = -18709.166
$ws.Range("N88").Value = -6947039
$ws.Range("H91").Value = 4637189.5
$ws.Range("I91").Value = 19115.166
$ws.Range("J91").Value = 6946227
$ws.Range("K91").Value = 19115.166
$ws.Range("L91").Value = 6946227
$ws.Range("M91").Value = -17711.166
$ws.Range("N91").Value = -6949035
$ws.Range("H110").Value = 2789.75
$ws.Range("I110").Value = 2789.75
$ws.Range("K110").Value = 2789.75
$ws.Range("M110").Value = -744.75
$ws.Range("H136").Value = 5689.5
$ws.Range("I136").Value = 4730.273
$ws.Range("K136").Value = 14190.819
$ws.Range("M136").Value = -11640.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6402.6
$ws.Range("I134").Value = 5062.4
$ws.Range("K134").Value = 15187.2
$ws.Range("M134").Value = -12652.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 15729.444
$ws.Range("J109").Value = 15729.444
$ws.Range("L109").Value = 15729.444
$ws.Range("N109").Value = -17809.444
$ws.Range("H134").Value = 5728.2144
$ws.Range("I134").Value = 3456.4285
$ws.Range("K134").Value = 10369.2855
$ws.Range("M134").Value = -7834.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 408
$ws.Range("J18").Value = 621
$ws.Range("L18").Value = 1863
$ws.Range("N18").Value = -2201
$ws.Range("H102").Value = 2800
$ws.Range("J102").Value = 2800
$ws.Range("L102").Value = 8400
$ws.Range("N102").Value = -13268
$ws.Range("H126").Value = 2499.25
$ws.Range("J126").Value = 2498
$ws.Range("L126").Value = 7494
$ws.Range("N126").Value = -17374
$ws.Range("H129").Value = 1289.4
$ws.Range("I129").Value = 891.25
$ws.Range("J129").Value = 1656.9231
$ws.Range("K129").Value = 2673.75
$ws.Range("L129").Value = 4970.7693
$ws.Range("M129").Value = 2326.25
$ws.Range("N129").Value = -14970.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40239.594
$ws.Range("I70").Value = 66663.89
$ws.Range("K70").Value = 66663.89
$ws.Range("M70").Value = -66393.89
$ws.Range("H73").Value = 40239.594
$ws.Range("I73").Value = 66663.89
$ws.Range("K73").Value = 66663.89
$ws.Range("M73").Value = -65727.89
$ws.Range("H80").Value = 15928835
$ws.Range("J80").Value = 66670444
$ws.Range("L80").Value = 66670444
$ws.Range("N80").Value = -66672440
$ws.Range("H83").Value = 15928835
$ws.Range("J83").Value = 66670444
$ws.Range("L83").Value = 333352220
$ws.Range("N83").Value = -333362204
$ws.Range("H132").Value = 5222.636
$ws.Range("I132").Value = 4241.2607
$ws.Range("J132").Value = 7479.8
$ws.Range("K132").Value = 12723.7821
$ws.Range("L132").Value = 22439.4
$ws.Range("M132").Value = -10193.7821
$ws.Range("N132").Value = -27499.4
$ws.Range("H136").Value = 18260.666
$ws.Range("J136").Value = 18260.666
$ws.Range("L136").Value = 54781.99800000001
$ws.Range("N136").Value = -59881.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3694.75
$ws.Range("I7").Value = 1823.6
$ws.Range("K7").Value = 1823.6
$ws.Range("M7").Value = -1711.6
$ws.Range("H46").Value = 8668.343999999999
$ws.Range("I46").Value = 4098.2856
$ws.Range("J46").Value = 9947.959999999999
$ws.Range("K46").Value = 4098.2856
$ws.Range("L46").Value = 9947.959999999999
$ws.Range("M46").Value = -3910.2856
$ws.Range("N46").Value = -10323.96
$ws.Range("H61").Value = 5893.875
$ws.Range("I61").Value = 3787.875
$ws.Range("J61").Value = 7999.875
$ws.Range("K61").Value = 3787.875
$ws.Range("L61").Value = 7999.875
$ws.Range("M61").Value = -3585.875
$ws.Range("N61").Value = -8403.875
$ws.Range("H93").Value = 8669717
$ws.Range("I93").Value = 2265.8696
$ws.Range("K93").Value = 2265.8696
$ws.Range("M93").Value = -1017.8696
$ws.Range("H113").Value = 5893.875
$ws.Range("I113").Value = 3787.875
$ws.Range("J113").Value = 7999.875
$ws.Range("K113").Value = 3787.875
$ws.Range("L113").Value = 7999.875
$ws.Range("M113").Value = -1617.875
$ws.Range("N113").Value = -12339.875
$ws.Range("H122").Value = 5105.516
$ws.Range("I122").Value = 4017.4285
$ws.Range("J122").Value = 7390.5
$ws.Range("K122").Value = 12052.2855
$ws.Range("L122").Value = 22171.5
$ws.Range("M122").Value = -9602.2855
$ws.Range("N122").Value = -27071.5
$ws.Range("H126").Value = 3694.75
$ws.Range("I126").Value = 1823.6
$ws.Range("K126").Value = 5470.799999999999
$ws.Range("M126").Value = -3000.799999999999
$ws.Range("H132").Value = 4662.1113
$ws.Range("I132").Value = 4623.7427
$ws.Range("K132").Value = 13871.2281
$ws.Range("M132").Value = -11341.2281

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64849
$ws.Range("J46").Value = 64849
$ws.Range("L46").Value = 64849
$ws.Range("N46").Value = -65311
$ws.Range("H51").Value = 23791.5
$ws.Range("I51").Value = 20999.777
$ws.Range("J51").Value = 32166.666
$ws.Range("K51").Value = 20999.777
$ws.Range("L51").Value = 32166.666
$ws.Range("M51").Value = -20489.777
$ws.Range("N51").Value = -33186.666
$ws.Range("H122").Value = 2810.6155
$ws.Range("I122").Value = 2235.8064
$ws.Range("K122").Value = 6707.4192
$ws.Range("M122").Value = -4257.4192
$ws.Range("H132").Value = 3310.318
$ws.Range("I132").Value = 1725.091
$ws.Range("J132").Value = 4895.5454
$ws.Range("K132").Value = 5175.272999999999
$ws.Range("L132").Value = 14686.6362
$ws.Range("M132").Value = -2645.272999999999
$ws.Range("N132").Value = -19746.6362
$ws.Range("H134").Value = 64849
$ws.Range("J134").Value = 64849
$ws.Range("L134").Value = 194547
$ws.Range("N134").Value = -199617
$ws.Range("H136").Value = 3617.359
$ws.Range("I136").Value = 2841.0667
$ws.Range("K136").Value = 8523.2001
$ws.Range("M136").Value = -5973.2001
